$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update swaps the two observation records currently stored in row 2
# and row 3 (the records themselves are unchanged but now listed in the
# opposite order), refreshes the Ost/Nord (Q/R) coordinates to new
# (rounded) values, and drops the Starttid/Sluttid (Z/AB) time values.

# --- Row 2 becomes the "Korallrot" record (previously in row 3) ---
$ws.Range("A2").Value = 111799311
$ws.Range("B2").Value = 96251
$ws.Range("E2").Value = 220093
$ws.Range("F2").Value = "Korallrot"
$ws.Range("G2").Value = "Corallorhiza trifida"
$ws.Range("H2").Value = "Châtel."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = ""
$ws.Range("Q2").Value = 513718
$ws.Range("R2").Value = 6704677
$ws.Range("S2").Value = 25
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 becomes the "Rödgul trumpetsvamp" record (previously in row 2) ---
$ws.Range("A3").Value = 111799186
$ws.Range("B3").Value = 89183
$ws.Range("E3").Value = 3215
$ws.Range("F3").Value = "Rödgul trumpetsvamp"
$ws.Range("G3").Value = "Craterellus lutescens"
$ws.Range("H3").Value = "(Fr.) Fr."
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "20"
$ws.Range("L3").ClearContents()
$ws.Range("Q3").Value = 513785
$ws.Range("R3").Value = 6704707
$ws.Range("S3").Value = 10
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
